$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '35.729.32'
$ws.Cells.Item(2, 5).Value = '  +3.58%  '
$ws.Cells.Item(3, 4).Value = '1.864.99'
$ws.Cells.Item(3, 5).Value = '  +2.97%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.37%  '
Set-TextValue $ws.Cells.Item(5, 4) '231.47'
$ws.Cells.Item(5, 5).Value = '  +2.63%  '
$ws.Cells.Item(6, 5).Value = '  +3.64%  '
$ws.Cells.Item(7, 5).Value = '  +0.35%  '
Set-TextValue $ws.Cells.Item(8, 4) '42.82'
$ws.Cells.Item(8, 5).Value = '  +11.90%  '
$ws.Cells.Item(9, 5).Value = '  +7.55%  '
$ws.Cells.Item(10, 5).Value = '  +3.41%  '
$ws.Cells.Item(11, 5).Value = '  +3.93%  '
$ws.Cells.Item(12, 4).Value = '2.136.04'
$ws.Cells.Item(12, 5).Value = '  +3.04%  '
Set-TextValue $ws.Cells.Item(13, 4) '11.73'
$ws.Cells.Item(13, 5).Value = '  +4.35%  '
$ws.Cells.Item(14, 4).Value = '1.871.74'
$ws.Cells.Item(14, 5).Value = '  +3.10%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.683'
$ws.Cells.Item(15, 5).Value = '  +7.90%  '
$ws.Cells.Item(16, 5).Value = '  +7.34%  '
$ws.Cells.Item(17, 4).Value = '35.744.80'
$ws.Cells.Item(17, 5).Value = '  +3.78%  '
Set-TextValue $ws.Cells.Item(18, 4) '70.66'
$ws.Cells.Item(18, 5).Value = '  +3.33%  '
Set-TextValue $ws.Cells.Item(19, 4) '249.43'
$ws.Cells.Item(19, 5).Value = '  +2.58%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0807'
$ws.Cells.Item(20, 5).Value = '  +4.56%  '
Set-TextValue $ws.Cells.Item(21, 4) '12.35'
$ws.Cells.Item(21, 5).Value = '  +10.10%  '
$ws.Cells.Item(22, 5).Value = '  +15.38%  '
$ws.Cells.Item(23, 5).Value = '  +0.28%  '
$ws.Cells.Item(24, 5).Value = '  +1.71%  '
Set-TextValue $ws.Cells.Item(25, 4) '170.94'
$ws.Cells.Item(25, 5).Value = '  +0.33%  '
Set-TextValue $ws.Cells.Item(26, 4) '8.03'
$ws.Cells.Item(26, 5).Value = '  +2.79%  '
Set-TextValue $ws.Cells.Item(27, 4) '17.92'
$ws.Cells.Item(27, 5).Value = '  +1.77%  '
$ws.Cells.Item(28, 5).Value = '  +1.92%  '
Set-TextValue $ws.Cells.Item(29, 4) '1.44'
$ws.Cells.Item(29, 5).Value = '  +16.78%  '
$ws.Cells.Item(30, 5).Value = '  +0.39%  '
$ws.Cells.Item(31, 4).Value = '3.324.08'
$ws.Cells.Item(31, 5).Value = '  +36.81%  '
Set-TextValue $ws.Cells.Item(35, 4) '1.91'
$ws.Cells.Item(35, 5).Value = '  +4.76%  '
Set-TextValue $ws.Cells.Item(36, 4) '100.68'
$ws.Cells.Item(36, 5).Value = '  +22.89%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.692'
$ws.Cells.Item(37, 5).Value = '  +7.61%  '
$ws.Cells.Item(39, 4).Value = '1.367.33'
$ws.Cells.Item(39, 5).Value = '  +0.44%  '
$ws.Cells.Item(40, 5).Value = '  +3.24%  '
$ws.Cells.Item(41, 5).Value = '  +5.28%  '
$ws.Cells.Item(42, 5).Value = '  +6.76%  '
Set-TextValue $ws.Cells.Item(43, 4) '15.03'
$ws.Cells.Item(43, 5).Value = '  +8.80%  '
Set-TextValue $ws.Cells.Item(44, 4) '1.26'
$ws.Cells.Item(44, 5).Value = '  +3.91%  '
Set-TextValue $ws.Cells.Item(45, 4) '2.47'
$ws.Cells.Item(45, 5).Value = '  +1.13%  '
Set-TextValue $ws.Cells.Item(46, 4) '2.83'
$ws.Cells.Item(46, 5).Value = '  +0.96%  '
Set-TextValue $ws.Cells.Item(47, 4) '6.30'
$ws.Cells.Item(47, 5).Value = '  +9.09%  '
Set-TextValue $ws.Cells.Item(48, 4) '0.0520'
$ws.Cells.Item(48, 5).Value = '  +2.29%  '
$ws.Cells.Item(49, 4).Value = '2.034.78'
$ws.Cells.Item(49, 5).Value = '  +3.07%  '
Set-TextValue $ws.Cells.Item(50, 4) '105.25'
$ws.Cells.Item(50, 5).Value = '  +2.72%  '
$ws.Cells.Item(51, 5).Value = '  +0.36%  '
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Cells.Item(32, 4) '4.11'
$ws.Cells.Item(32, 5).Value = '  +6.43%  '
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Cells.Item(33, 4) '0.0549'
$ws.Cells.Item(33, 5).Value = '  +6.11%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(34, 4) '3.97'
$ws.Cells.Item(34, 5).Value = '  +4.61%  '

Write-Host "Applied all changes"
